# Add I0 and IF columns (I and J) to the worksheet, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, styled like the existing header row (s="1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-52 for columns I (I0) and J (IF)
$data = @(
    @{Row=2; I=8; J=8},
    @{Row=3; I=7; J=7},
    @{Row=4; I=7; J=7},
    @{Row=5; I=7; J=7},
    @{Row=6; I=7; J=7},
    @{Row=7; I=7; J=7},
    @{Row=8; I=7; J=7},
    @{Row=9; I=8; J=8},
    @{Row=10; I=1; J=1},
    @{Row=11; I=8; J=8},
    @{Row=12; I=1; J=2},
    @{Row=13; I=7; J=7},
    @{Row=14; I=6; J=6},
    @{Row=15; I=9; J=9},
    @{Row=16; I=7; J=8},
    @{Row=17; I=6; J=6},
    @{Row=18; I=6; J=6},
    @{Row=19; I=6; J=7},
    @{Row=20; I=8; J=8},
    @{Row=21; I=6; J=7},
    @{Row=22; I=7; J=7},
    @{Row=23; I=8; J=8},
    @{Row=24; I=8; J=8},
    @{Row=25; I=7; J=7},
    @{Row=26; I=6; J=6},
    @{Row=27; I=7; J=7},
    @{Row=28; I=7; J=7},
    @{Row=29; I=6; J=6},
    @{Row=30; I=6; J=7},
    @{Row=31; I=8; J=9},
    @{Row=32; I=7; J=7},
    @{Row=33; I=9; J=9},
    @{Row=34; I=8; J=8},
    @{Row=35; I=7; J=7},
    @{Row=36; I=7; J=8},
    @{Row=37; I=6; J=7},
    @{Row=38; I=7; J=7},
    @{Row=39; I=7; J=7},
    @{Row=40; I=7; J=7},
    @{Row=41; I=8; J=8},
    @{Row=42; I=8; J=8},
    @{Row=43; I=6; J=7},
    @{Row=44; I=8; J=8},
    @{Row=45; I=5; J=5},
    @{Row=46; I=8; J=8},
    @{Row=47; I=4; J=5},
    @{Row=48; I=8; J=8},
    @{Row=49; I=7; J=7},
    @{Row=50; I=8; J=8},
    @{Row=51; I=5; J=5},
    @{Row=52; I=8; J=8}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}
